{"js": "const body = context.document.body;\n\n// Update the title line (date) at the top of the document.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\nif (titlePara.text.trim() === \"2023-09-29 Friday\") {\n  titlePara.getRange().insertText(\"2023-09-30 Saturday\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Update each arithmetic-problem cell in the practice table (20 rows x 5 cols).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"50-9=\", \"72-29=\", \"46+36=\", \"81-33=\", \"58+23=\"],\n  [\"82-18=\", \"27+64=\", \"88+5=\", \"70-62=\", \"84-48=\"],\n  [\"51-37=\", \"18+59=\", \"50-9=\", \"58-19=\", \"75-8=\"],\n  [\"96-18=\", \"62-49=\", \"40-19=\", \"49+19=\", \"25-7=\"],\n  [\"63+28=\", \"84-37=\", \"3+68=\", \"5+87=\", \"44-9=\"],\n  [\"7+25=\", \"6+29=\", \"23+59=\", \"97-19=\", \"98-29=\"],\n  [\"63-59=\", \"48+39=\", \"55+17=\", \"82-46=\", \"64-49=\"],\n  [\"44+18=\", \"85-27=\", \"86+8=\", \"20-4=\", \"69+4=\"],\n  [\"58+7=\", \"76-7=\", \"85-6=\", \"46+9=\", \"93-25=\"],\n  [\"63-19=\", \"64+27=\", \"76-18=\", \"74+19=\", \"70-55=\"],\n  [\"84-15=\", \"19+64=\", \"6+29=\", \"74+19=\", \"36-28=\"],\n  [\"26+18=\", \"30-12=\", \"28+9=\", \"70-54=\", \"47+5=\"],\n  [\"17+38=\", \"95-76=\", \"46+47=\", \"56+17=\", \"80-59=\"],\n  [\"33-29=\", \"53-4=\", \"55-38=\", \"70-52=\", \"74-28=\"],\n  [\"59+38=\", \"44+8=\", \"6+8=\", \"8+37=\", \"92-46=\"],\n  [\"73-64=\", \"96-27=\", \"27-8=\", \"90-79=\", \"47+48=\"],\n  [\"16+7=\", \"73-59=\", \"75-16=\", \"39+44=\", \"48+23=\"],\n  [\"39+22=\", \"88+5=\", \"58+8=\", \"78+18=\", \"92-35=\"],\n  [\"27+39=\", \"54+27=\", \"94-17=\", \"86-17=\", \"89+6=\"],\n  [\"38+29=\", \"94-76=\", \"9+17=\", \"32-17=\", \"81-28=\"],\n];\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(r, c);\n    cell.body.getRange().insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n# Update the title line (date) at the top of the document.\n$titlePara = $d.Paragraphs.Item(1)\n$titleOld = \"2023-09-29 Friday\"\n$titleNew = \"2023-09-30 Saturday\"\nif ($titlePara.Range.Text.TrimEnd([char]13,[char]7) -eq $titleOld) {\n    $titlePara.Range.Text = $titleNew\n}\n\n# Update each arithmetic-problem cell in the practice table (20 rows x 5 cols).\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"50-9=\", \"72-29=\", \"46+36=\", \"81-33=\", \"58+23=\"),\n    @(\"82-18=\", \"27+64=\", \"88+5=\", \"70-62=\", \"84-48=\"),\n    @(\"51-37=\", \"18+59=\", \"50-9=\", \"58-19=\", \"75-8=\"),\n    @(\"96-18=\", \"62-49=\", \"40-19=\", \"49+19=\", \"25-7=\"),\n    @(\"63+28=\", \"84-37=\", \"3+68=\", \"5+87=\", \"44-9=\"),\n    @(\"7+25=\", \"6+29=\", \"23+59=\", \"97-19=\", \"98-29=\"),\n    @(\"63-59=\", \"48+39=\", \"55+17=\", \"82-46=\", \"64-49=\"),\n    @(\"44+18=\", \"85-27=\", \"86+8=\", \"20-4=\", \"69+4=\"),\n    @(\"58+7=\", \"76-7=\", \"85-6=\", \"46+9=\", \"93-25=\"),\n    @(\"63-19=\", \"64+27=\", \"76-18=\", \"74+19=\", \"70-55=\"),\n    @(\"84-15=\", \"19+64=\", \"6+29=\", \"74+19=\", \"36-28=\"),\n    @(\"26+18=\", \"30-12=\", \"28+9=\", \"70-54=\", \"47+5=\"),\n    @(\"17+38=\", \"95-76=\", \"46+47=\", \"56+17=\", \"80-59=\"),\n    @(\"33-29=\", \"53-4=\", \"55-38=\", \"70-52=\", \"74-28=\"),\n    @(\"59+38=\", \"44+8=\", \"6+8=\", \"8+37=\", \"92-46=\"),\n    @(\"73-64=\", \"96-27=\", \"27-8=\", \"90-79=\", \"47+48=\"),\n    @(\"16+7=\", \"73-59=\", \"75-16=\", \"39+44=\", \"48+23=\"),\n    @(\"39+22=\", \"88+5=\", \"58+8=\", \"78+18=\", \"92-35=\"),\n    @(\"27+39=\", \"54+27=\", \"94-17=\", \"86-17=\", \"89+6=\"),\n    @(\"38+29=\", \"94-76=\", \"9+17=\", \"32-17=\", \"81-28=\"),\n)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $row[$c - 1]\n    }\n}"}
